# The workbook contains a single sheet with a table of Zanahoria (carrot)
# price observations. A new daily observation needs to be inserted as a new
# row 284 (pushing the existing rows 284-399 down to 285-400).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 284; existing rows 284.. shift down by one.
$ws.Rows(284).Insert()

# Populate the newly inserted row 284 with the new observation.
$ws.Range("A284").Value = 3
$ws.Range("B284").Value = 'Femacal de La Calera'
$ws.Range("C284").Value = 'Coquimbo'
$ws.Range("D284").Value = 44784
$ws.Range("E284").Value = 5
$ws.Range("F284").Value = 100114013
$ws.Range("G284").Value = 'Zanahoria'
$ws.Range("H284").Value = 'Sin especificar'
$ws.Range("I284").Value = 'Primera'
$ws.Range("J284").Value = 510
$ws.Range("K284").Value = 10000
$ws.Range("L284").Value = 11000
$ws.Range("M284").Value = 10510
$ws.Range("N284").Value = '$/saco 20 kilos'
$ws.Range("O284").Value = 'Provincia de Quillota'
$ws.Range("P284").Value = 526
$ws.Range("Q284").Value = 20
$ws.Range("R284").Value = 'Hortaliza'
